# Clear out the now-unused per-status breakdown columns (AA:AL) and the
# DIFFERENCE column (AN) for data rows 2-5 on Sheet1, keeping the
# PREVIOUS ACCOMPLISHMENT column (AM) intact, per the updated
# status/accomplishment source file (as of May).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear columns AA through AL (12 columns) for rows 2-5
$ws.Range("AA2:AL5").ClearContents()

# Clear column AN (DIFFERENCE) for rows 2-5, leave AM (PREVIOUS ACCOMPLISHMENT) untouched
$ws.Range("AN2:AN5").ClearContents()
